# Fuel Consumption Mapping - make the data tables more readable.
#
# The sheet holds three little x/y datasets (two in columns P:Q, one in
# columns E:F) that feed three scatter charts. This adds a descriptive
# "Dataset N" banner above each table and a friendlier "Speed (mph)" /
# "Desired MPG" sub-heading above the existing terse "x"/"y" headers, then
# centers all of the header/data cells and gives the charts real titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert one blank row above the first table -------------
# Before: table headers started on row 6 (E/F) & row 6 (P/Q) & row 25 (P/Q).
# After inserting a row at 5, those become row 7 / row 7 / row 26, and row 5
# (brand new) + row 6 (previously always-empty) are free to hold the new
# banner / sub-heading text.
[void]$ws.Rows("5:5").Insert()

# --- 2. New sub-heading row (row 6) above each "x"/"y" header -------------
$ws.Range("E6").Value = "Speed (mph)"
$ws.Range("F6").Value = "Desired MPG"
$ws.Range("P6").Value = "Speed (mph)"
$ws.Range("Q6").Value = "Desired MPG"

# --- 3. New "Dataset N" banner row for each table -------------------------
$ws.Range("P24").Value = "Dataset 3"
$ws.Range("P5").Value = "Dataset 2"
$ws.Range("E5").Value = "Dataset 1"

# --- 4. Merge the banner cells across both columns of their table --------
[void]$ws.Range("P24:Q24").Merge()
[void]$ws.Range("P5:Q5").Merge()
[void]$ws.Range("E5:F5").Merge()

# --- 5. Center the banner rows ---------------------------------------------
$ws.Range("E5:F5").HorizontalAlignment = -4108
$ws.Range("P5:Q5").HorizontalAlignment = -4108
$ws.Range("P24:Q24").HorizontalAlignment = -4108

# --- 6. Center every header + data cell in the three tables ---------------
# Table 1 (E/F): header row 7, data rows 8-16
$ws.Range("E7:F16").HorizontalAlignment = -4108
# Table 2 (P/Q): header row 7, data rows 8-21
$ws.Range("P7:Q21").HorizontalAlignment = -4108
# Table 3 (P/Q): header row 26, data rows 27-38
$ws.Range("P26:Q38").HorizontalAlignment = -4108

# --- 7. Let the widened headers dictate sensible column widths ------------
[void]$ws.Range("E1:F1").EntireColumn.AutoFit()
[void]$ws.Range("P1:Q1").EntireColumn.AutoFit()

# --- 8. Give each chart a real title instead of a blank placeholder -------
for ($i = 1; $i -le $ws.ChartObjects().Count; $i++) {
    $co = $ws.ChartObjects($i)
    $co.Chart.HasTitle = $true
    $co.Chart.ChartTitle.Text = "Fuel Economy vs. Speed"
}

# --- 9. Restore a sensible selection/view ----------------------------------
[void]$ws.Range("R15").Select()
